$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-111). Bump it from 45203 (2023-10-04) to 45204 (2023-10-05) for
# every row that currently has the old value, leaving anything else intact.
for ($r = 2; $r -le 111; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
